$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: refreshed timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 08:24"

# Row 4: Estados Unidos
$ws.Range("B4").Value = 7406729
$ws.Range("C4").Value = 583
$ws.Range("D4").Value = 4649521
$ws.Range("E4").Value = 2546411
$ws.Range("G4").Value = 12
$ws.Range("H4").Value = 210797

# Row 60: Uzbekistan
$ws.Range("B60").Value = 56519
$ws.Range("C60").Value = 165
$ws.Range("E60").Value = 3081
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 468

# Row 77: El Salvador
$ws.Range("E77").Value = 4453
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 843

# Row 80: Australia
$ws.Range("D80").Value = 24754
$ws.Range("E80").Value = 1438

# Rows 215-216: Montserrat now listed before Islas Malvinas
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
